$wb = $excel.ActiveWorkbook

# 1. Delete the two extra (empty) sheets, keeping only the first one.
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Sheet3").Delete()
$wb.Worksheets.Item("Sheet2").Delete()

# 2. Rename the remaining sheet.
$ws = $wb.Worksheets.Item(1)
$ws.Name = "OSM Order"

# 3. Update the bill's customer / address / phone details.
$ws.Range("A4").Value = "Customer Name: User 1"
$ws.Range("A5").Value = "Address: 228/21 abc, xyz"
$ws.Range("A6").Value = "Phone: 12312321312"

# 4. Update the line item.
$ws.Range("B9").Value = "Toy 7"
$ws.Range("C9").Value = "1"
$ws.Range("D9").Value = "150,000"
$ws.Range("E9").Value = "150,000"

# 5. Update the date field near the signature block.
$ws.Range("C28").Value = "15, 5, 2019"

# 6. Update the active selection on the sheet.
$ws.Range("G14").Select()
